# Append a new day (2020-04-15 / serial 43936) of age-bucket rows to the
# TN_AgeDaily sheet, continuing the existing fill-down formulas in
# columns D, E, F, H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 281
$newDate = $ws.Cells.Item($lastRow, 1).Value2 + 1

# Age bucket labels (column B) in the order they repeat for every date block.
$labels = @("0-10", "11-20", "21-30", "31-40", "41-50", "51-60", "61-70", "71-80", "80+", "Pending")

# TOT_CASE_COUNT (C) and DEATHS_TOT (G) values for the new date, one per label.
$totCase = @(70, 272, 1189, 958, 969, 1076, 731, 353, 197, 264)
$deathsTot = @(1, 0, 1, 1, 7, 14, 34, 33, 44, 0)

$newFirstRow = $lastRow + 1
$newLastRow = $lastRow + $labels.Length

# Carry the existing row's number formats (date / integer / percent) down
# onto the new block before filling in values and formulas.
$ws.Range("A$lastRow`:H$lastRow").Copy()
$ws.Range("A$newFirstRow`:H$newLastRow").PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $labels.Length; $i++) {
    $r = $newFirstRow + $i

    $ws.Cells.Item($r, 1).Value2 = $newDate
    $ws.Cells.Item($r, 2).Value2 = $labels[$i]
    $ws.Cells.Item($r, 3).Value2 = $totCase[$i]

    $ws.Cells.Item($r, 4).Formula = "=C$r/SUMIF(A:A,A$r,C:C)"
    $ws.Cells.Item($r, 5).Formula = "=C$r-SUMIFS(C:C,A:A,A$r-1,B:B,B$r)"
    $ws.Cells.Item($r, 6).Formula = "=E$r/SUMIF(A:A,A$r,E:E)"

    $ws.Cells.Item($r, 7).Value2 = $deathsTot[$i]

    $ws.Cells.Item($r, 8).Formula = "=G$r-SUMIFS(G:G,A:A,A$r-1,B:B,B$r)"
}

$ws.Range("H287:H290").Select()
